$d = $word.ActiveDocument

# 1) "In contrast to Americans are ... Polish are distressful." ->
#    "Whilst Americans are ... Poles are distressful."
$d.Content.Find.Execute(
    "In contrast to Americans are cheerful and optimistic, Polish are distressful.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Whilst Americans are cheerful and optimistic, Poles are distressful.", 2) | Out-Null

# 2) "Majority of American are easy-going and sociable, unlike majority of
#     Poles is reserved and often antisocial." ->
#    "Majority of Americans are easy-going and sociable, unlike majority of
#     Poles, who are reserved and often antisocial."
$d.Content.Find.Execute(
    "Majority of American are easy-going and sociable, unlike majority of Poles is reserved and often antisocial.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Majority of Americans are easy-going and sociable, unlike majority of Poles, who are reserved and often antisocial.", 2) | Out-Null

# 3) "Stereotypical American is fat and like eating hamburgers, unlike Poles
#     likes pork chop and potatoes." ->
#    "Stereotypical American is fat and like eating hamburgers, unlike
#     cut-and-dried Pole, who likes pork chop and potatoes."
$d.Content.Find.Execute(
    "Stereotypical American is fat and like eating hamburgers, unlike Poles likes pork chop and potatoes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Stereotypical American is fat and like eating hamburgers, unlike cut-and-dried Pole, who likes pork chop and potatoes.", 2) | Out-Null

# 4) The "_GoBack" bookmark used to sit right after the second picture
#    (end of that paragraph); it now belongs in the middle of the rewritten
#    last paragraph, between "cut-and-" and "dried". Move it there.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$lastPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$findRange = $lastPara.Range.Duplicate
$findRange.Find.Execute("dried", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$bmRange = $d.Range($findRange.Start, $findRange.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
